# EMEP_NFR09_scaling_mapping.xlsx - add a new scaling rule row to the "year" sheet
# and nudge the saved window position, per the commit "Updates to scaling
# mappings and scaling rules".

$wb = $excel.ActiveWorkbook

# Record the window position that will be stored with the workbook view.
$win = $wb.Windows.Item(1)
$win.Left = 6680
$win.Top = 2280

$ws = $wb.Worksheets.Item("year")

# New row 4: lux / RoadRail scaling rule, 1990-2020, with a comment.
$ws.Range("A4").Value = "lux"
$ws.Range("B4").Value = "RoadRail"
$ws.Range("C4").Value = "NA"
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = 1990
$ws.Range("G4").Value = 2020
$ws.Range("H4").Value = "NA"
$ws.Range("I4").Value = "Avoid imlied Nox EF dip 1986-1989"

# Match the author's final selection (whole new row highlighted).
$ws.Range("A4:XFD4").Select()
